$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing row 34 ---
# B34 was a time value (0.770833...) with style 5; it becomes a text string, keeping style 5
$ws.Range("B34").Value = "18.30-20.30"

# F34 (new content - "META" column)
$ws.Range("F34").Value = "Paljon taas teoriaa tankattu, mutta ehkä se tästä sitte pikkuhiljaa maturoituu"
$ws.Range("F34").WrapText = $true

# D34 (new content - "Oppimisen laatu" column)
$ws.Range("D34").Value = "Vaatii uudelleenlukemista, mutta virkeänä ja valppaana keskittyneesti eteenpäin, pieniä sivupolkuja unohtamatta."
$ws.Range("D34").WrapText = $true

# --- Add new row 35 ---
$ws.Range("A35").Value = "26 marras"

$ws.Range("C35").Value = "Siirtyminen törmäyksen havaitsemisesta kontaktien aiheuttamiin voimiin"

# C34 gets new, extended text (set after row 35's strings, matching shared-string order)
$ws.Range("C34").Value = "Erottavan hypertason teoreema, kahden monikulmion leikkaustarkastelu (box-box intersection), monikulmion laajennusta, s. 279-290"

# G34 (hours worked)
$ws.Range("G34").Value = 2

# Row 34 height grows to fit the new wrapped content
$ws.Rows.Item(34).RowHeight = 72.5

$ws.Range("B35").Value = 0.39583333333333331
$ws.Range("B35").NumberFormat = "h:mm"

$ws.Rows.Item(35).RowHeight = 43.5
